# Update crypto price (column D) and 1h volume change (column E) values
# to reflect the refreshed data pulled by the scheduled GitHub Actions job.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> @{ D = "<new price>"; E = "<new pct>" }
$updates = @{
    2  = @{ D = "301.10";       E = "0.67%" }
    3  = @{ D = "32.16";        E = "1.53%" }
    4  = @{ D = "4.979";        E = "-3.14%" }
    5  = @{              E = "-2.24%" }
    6  = @{ D = "2.104";        E = "-16.10%" }
    7  = @{ D = "7.798";        E = "-0.07%" }
    8  = @{ D = "3.845";        E = "-1.66%" }
    9  = @{ D = "0.9258";       E = "-0.23%" }
    10 = @{ D = "0.1746";       E = "-0.63%" }
    11 = @{ D = "0.07972";      E = "7.25%" }
    12 = @{ D = "0.08679";      E = "-1.79%" }
    13 = @{ D = "0.03106";      E = "3.54%" }
    14 = @{ D = "0.1003";       E = "0.26%" }
    15 = @{ D = "0.001525";     E = "1.12%" }
    16 = @{ D = "0.005865";     E = "-1.45%" }
    17 = @{              E = "2,100.91%" }
    18 = @{ D = "3.460";        E = "-2.73%" }
    19 = @{              E = "-1.44%" }
    20 = @{ D = "0.3287";       E = "0.48%" }
    21 = @{ D = "0.1309";       E = "-2.19%" }
    22 = @{ D = "4.314";        E = "3.44%" }
    23 = @{ D = "0.1792";       E = "6.68%" }
    24 = @{ D = "0.04605";      E = "-0.36%" }
    25 = @{              E = "-0.12%" }
    26 = @{ D = "0.004441";     E = "-1.92%" }
    27 = @{ D = "0.0001251";    E = "4.18%" }
    39 = @{ D = "0.01714";      E = "-2.22%" }
    40 = @{ D = "0.04779";      E = "3.97%" }
    41 = @{ D = "0.007500";     E = "8.35%" }
    42 = @{ D = "0.1359";       E = "-0.87%" }
    43 = @{ D = "0.002351";     E = "7.32%" }
    44 = @{ D = "0.01125";      E = "4.45%" }
    45 = @{ D = "0.00005997";   E = "-2.09%" }
    46 = @{ D = "0.00000000750"; E = "0.12%" }
    47 = @{ D = "0.003391";     E = "-59.61%" }
    49 = @{ D = "0.00002101";   E = "0.12%" }
    50 = @{ D = "0.0002001";    E = "0.12%" }
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    if ($vals.ContainsKey("D")) {
        $cell = $ws.Cells.Item($row, 4)
        # Force text storage so values like "301.10" / "0.07500" keep their
        # exact original formatting instead of being coerced to numbers.
        $cell.NumberFormat = "@"
        $cell.Value = $vals["D"]
    }
    if ($vals.ContainsKey("E")) {
        $cell = $ws.Cells.Item($row, 5)
        $cell.NumberFormat = "@"
        $cell.Value = $vals["E"]
    }
}
